$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns D, L, M, N, O, P, Q, S).
# Derived from the target diff: the weekly price records were reshuffled
# across rows (dates/quality/volume/prices/unit/price-per-kg moved between
# rows while the market/product metadata columns stayed put).
$rows = @(
    @{ Row = 2;  D = 45044; L = "Primera"; M = 100; N = 17000; O = 18000; P = 17500; Q = "`$/caja 18 kilos";        S = 972  }
    @{ Row = 3;  D = 45030; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; S = 861  }
    @{ Row = 4;  D = 44819; L = "Primera"; M = 100; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos granel"; S = 1417 }
    @{ Row = 5;  D = 45014; L = "Primera"; M = 50;  N = 13000; O = 14000; P = 13600; Q = "`$/caja 18 kilos";        S = 756  }
    @{ Row = 6;  D = 45014; L = "Segunda"; M = 20;  N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos";        S = 556  }
    @{ Row = 7;  D = 44316; L = "Primera"; M = 50;  N = 20000; O = 20000; P = 20000; Q = "`$/caja 18 kilos";        S = 1111 }
    @{ Row = 8;  D = 44699; L = "Primera"; M = 100; N = 20000; O = 22000; P = 21000; Q = "`$/caja 18 kilos";        S = 1167 }
    @{ Row = 9;  D = 44699; L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos";        S = 1000 }
    @{ Row = 10; D = 45002; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "`$/caja 18 kilos";        S = 694  }
    @{ Row = 11; D = 44516; L = "Primera"; M = 100; N = 33000; O = 34000; P = 33500; Q = "`$/caja 18 kilos";        S = 1861 }
    @{ Row = 12; D = 45084; L = "Primera"; M = 100; N = 20000; O = 21000; P = 20500; Q = "`$/caja 18 kilos granel"; S = 1139 }
    @{ Row = 13; D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos";        S = 806  }
    @{ Row = 14; D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";        S = 667  }
    @{ Row = 15; D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "`$/caja 18 kilos";        S = 1028 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Precio mínimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Precio máximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Unidad de comercialización
    $ws.Cells.Item($row, 19).Value = $r.S   # S: Precio $/Kg
}
